# Auto-generated cell updates matching the target OOXML diff.
# Values are written via Range.Value; $null clears a cell entirely
# (matches cells removed from the XML in the diff).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 28
$ws.Range("H28").Value = 9049
$ws.Range("J28").Value = 8098
$ws.Range("L28").Value = 8098
$ws.Range("N28").Value = -9068
# row 62
$ws.Range("H62").Value = 18522836
$ws.Range("I62").Value = 22226802
$ws.Range("K62").Value = 22226802
$ws.Range("M62").Value = -22226178
# row 65
$ws.Range("H65").Value = 18522836
$ws.Range("I65").Value = 22226802
$ws.Range("K65").Value = 111134010
$ws.Range("M65").Value = -111130890
# row 69
$ws.Range("H69").Value = 2660
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 3490
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 10470
$ws.Range("M69").Value = -2126
$ws.Range("N69").Value = -12218
# row 72
$ws.Range("H72").Value = 2660
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 3490
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 31410
$ws.Range("M72").Value = -4632
$ws.Range("N72").Value = -40146
# row 76
$ws.Range("H76").Value = 8333.333000000001
$ws.Range("I76").Value = 10000
$ws.Range("J76").Value = 7500
$ws.Range("K76").Value = 10000
$ws.Range("L76").Value = 7500
$ws.Range("M76").Value = -9685
$ws.Range("N76").Value = -8130
# row 79
$ws.Range("H79").Value = 8333.333000000001
$ws.Range("I79").Value = 10000
$ws.Range("J79").Value = 7500
$ws.Range("K79").Value = 10000
$ws.Range("L79").Value = 7500
$ws.Range("M79").Value = -8908
$ws.Range("N79").Value = -9684
# row 107
$ws.Range("H107").Value = 9047.799999999999
$ws.Range("I107").Value = 8809.75
$ws.Range("K107").Value = 8809.75
$ws.Range("M107").Value = -6889.75
# row 111
$ws.Range("H111").Value = 1577.7778
$ws.Range("I111").Value = 1712.5
$ws.Range("J111").Value = 500
$ws.Range("K111").Value = 5137.5
$ws.Range("L111").Value = 1500
$ws.Range("M111").Value = -2070.5
$ws.Range("N111").Value = -7634
# row 135
$ws.Range("H135").Value = 32259512
$ws.Range("I135").Value = 604.95654
$ws.Range("J135").Value = 125003864
$ws.Range("K135").Value = 5444.60886
$ws.Range("L135").Value = 1125034776
$ws.Range("M135").Value = -2909.60886
$ws.Range("N135").Value = -1125039846
# row 137
$ws.Range("H137").Value = 2274.6785
$ws.Range("I137").Value = 2007.1052
$ws.Range("J137").Value = 2839.5557
$ws.Range("K137").Value = 6021.3156
$ws.Range("L137").Value = 8518.667099999999
$ws.Range("M137").Value = -3471.3156
$ws.Range("N137").Value = -13618.6671
# row 141
$ws.Range("H141").Value = 608.9583
$ws.Range("I141").Value = 608.9583
$ws.Range("K141").Value = 1826.8749
$ws.Range("M141").Value = 3353.1251

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 10442.551
$ws.Range("I32").Value = 7685.778
$ws.Range("J32").Value = 18076.691
$ws.Range("K32").Value = 7685.778
$ws.Range("L32").Value = 18076.691
$ws.Range("M32").Value = -7398.778
$ws.Range("N32").Value = -18650.691
# row 63
$ws.Range("H63").Value = 2082.2856
$ws.Range("I63").Value = 1939.9524
$ws.Range("J63").Value = 2509.2856
$ws.Range("K63").Value = 1939.9524
$ws.Range("L63").Value = 2509.2856
$ws.Range("M63").Value = -1253.9524
$ws.Range("N63").Value = -3881.2856
# row 66
$ws.Range("H66").Value = 2082.2856
$ws.Range("I66").Value = 1939.9524
$ws.Range("J66").Value = 2509.2856
$ws.Range("K66").Value = 9699.761999999999
$ws.Range("L66").Value = 12546.428
$ws.Range("M66").Value = -6267.761999999999
$ws.Range("N66").Value = -19410.428
# row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
# row 110
$ws.Range("H110").Value = 994.75
$ws.Range("I110").Value = 993
$ws.Range("K110").Value = 993
$ws.Range("M110").Value = 1052
# row 122
$ws.Range("H122").Value = 2771.8
$ws.Range("I122").Value = 2660
$ws.Range("J122").Value = 2939.5
$ws.Range("K122").Value = 7980
$ws.Range("L122").Value = 8818.5
$ws.Range("M122").Value = -5530
$ws.Range("N122").Value = -13718.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 99
$ws.Range("H99").Value = 58824690
$ws.Range("J99").Value = 1412.2
$ws.Range("L99").Value = 1412.2
$ws.Range("N99").Value = -4408.2
# row 105
$ws.Range("H105").Value = 142858190
$ws.Range("I105").Value = 166667570
$ws.Range("K105").Value = 166667570
$ws.Range("M105").Value = -166665823
# row 134
$ws.Range("H134").Value = 5846.909
$ws.Range("I134").Value = 1012.0769
$ws.Range("K134").Value = 3036.2307
$ws.Range("M134").Value = -501.2307000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 62
$ws.Range("H62").Value = 15386810
$ws.Range("I62").Value = 2338.6365
$ws.Range("K62").Value = 2338.6365
$ws.Range("M62").Value = -1714.6365
# row 65
$ws.Range("H65").Value = 15386810
$ws.Range("I65").Value = 2338.6365
$ws.Range("K65").Value = 11693.1825
$ws.Range("M65").Value = -8573.182500000001
# row 107
$ws.Range("H107").Value = 1202.5
$ws.Range("I107").Value = 598.75
$ws.Range("K107").Value = 598.75
$ws.Range("M107").Value = 1321.25
# row 132
$ws.Range("H132").Value = 1744.35
$ws.Range("I132").Value = 1445.88
$ws.Range("J132").Value = 2241.8
$ws.Range("K132").Value = 4337.64
$ws.Range("L132").Value = 6725.400000000001
$ws.Range("M132").Value = -1807.64
$ws.Range("N132").Value = -11785.4
# row 134
$ws.Range("H134").Value = 14287324
$ws.Range("I134").Value = 1638.92
$ws.Range("J134").Value = 50001536
$ws.Range("K134").Value = 4916.76
$ws.Range("L134").Value = 150004608
$ws.Range("M134").Value = -2381.76
$ws.Range("N134").Value = -150009678

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 63
$ws.Range("H63").Value = 15000.1
$ws.Range("J63").Value = 15000.1
$ws.Range("L63").Value = 15000.1
$ws.Range("N63").Value = -16372.1
# row 64
$ws.Range("H64").Value = 16200
$ws.Range("J64").Value = 16200
$ws.Range("L64").Value = 16200
$ws.Range("N64").Value = -16696
# row 66
$ws.Range("H66").Value = 15000.1
$ws.Range("J66").Value = 15000.1
$ws.Range("L66").Value = 45000.3
$ws.Range("N66").Value = -51864.3
# row 67
$ws.Range("H67").Value = 16200
$ws.Range("J67").Value = 16200
$ws.Range("L67").Value = 16200
$ws.Range("N67").Value = -17916
# row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = $null
# row 70
$ws.Range("H70").Value = 75003320
$ws.Range("I70").Value = 125002504
$ws.Range("J70").Value = 50003730
$ws.Range("K70").Value = 125002504
$ws.Range("L70").Value = 50003730
$ws.Range("M70").Value = -125002234
$ws.Range("N70").Value = -50004270
# row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = $null
# row 73
$ws.Range("H73").Value = 75003320
$ws.Range("I73").Value = 125002504
$ws.Range("J73").Value = 50003730
$ws.Range("K73").Value = 125002504
$ws.Range("L73").Value = 50003730
$ws.Range("M73").Value = -125001568
$ws.Range("N73").Value = -50005602
# row 74
$ws.Range("H74").Value = 56850
$ws.Range("J74").Value = 56850
$ws.Range("L74").Value = 56850
$ws.Range("N74").Value = -58722
# row 77
$ws.Range("H77").Value = 56850
$ws.Range("J77").Value = 56850
$ws.Range("L77").Value = 170550
$ws.Range("N77").Value = -179910
# row 80
$ws.Range("H80").Value = 4285.5
$ws.Range("I80").Value = 3005
$ws.Range("J80").Value = 4427.778
$ws.Range("K80").Value = 3005
$ws.Range("L80").Value = 4427.778
$ws.Range("M80").Value = -2007
$ws.Range("N80").Value = -6423.778
# row 82
$ws.Range("H82").Value = 16500
$ws.Range("J82").Value = 16500
$ws.Range("L82").Value = 16500
$ws.Range("N82").Value = -17266
# row 83
$ws.Range("H83").Value = 4285.5
$ws.Range("I83").Value = 3005
$ws.Range("J83").Value = 4427.778
$ws.Range("K83").Value = 15025
$ws.Range("L83").Value = 22138.89
$ws.Range("M83").Value = -10033
$ws.Range("N83").Value = -32122.89
# row 85
$ws.Range("H85").Value = 16500
$ws.Range("J85").Value = 16500
$ws.Range("L85").Value = 16500
$ws.Range("N85").Value = -19152
# row 86
$ws.Range("H86").Value = 39000
$ws.Range("J86").Value = 39000
$ws.Range("L86").Value = 39000
$ws.Range("N86").Value = -41372
# row 88
$ws.Range("H88").Value = 40000
$ws.Range("J88").Value = 40000
$ws.Range("L88").Value = 40000
$ws.Range("N88").Value = -40902
# row 89
$ws.Range("H89").Value = 39000
$ws.Range("J89").Value = 39000
$ws.Range("L89").Value = 117000
$ws.Range("N89").Value = -128856
# row 91
$ws.Range("H91").Value = 40000
$ws.Range("J91").Value = 40000
$ws.Range("L91").Value = 40000
$ws.Range("N91").Value = -43120
# row 97
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").Value = $null
# row 113
$ws.Range("H113").Value = 1575
$ws.Range("I113").Value = 1350
$ws.Range("K113").Value = 1350
$ws.Range("M113").Value = 820
# row 132
$ws.Range("H132").Value = 6396.0347
$ws.Range("I132").Value = 9048.799999999999
$ws.Range("J132").Value = 3553.7856
$ws.Range("K132").Value = 27146.4
$ws.Range("L132").Value = 10661.3568
$ws.Range("M132").Value = -24616.4
$ws.Range("N132").Value = -15721.3568

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 93
$ws.Range("H93").Value = 1022.8823
$ws.Range("I93").Value = 962.3333
$ws.Range("J93").Value = 1477
$ws.Range("K93").Value = 962.3333
$ws.Range("L93").Value = 1477
$ws.Range("M93").Value = 285.6667
$ws.Range("N93").Value = -3973
# row 132
$ws.Range("H132").Value = 2330.6858
$ws.Range("I132").Value = 1896
$ws.Range("K132").Value = 5688
$ws.Range("M132").Value = -3158

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 3538.3408
$ws.Range("I132").Value = 3852.4707
$ws.Range("J132").Value = 2470.3
$ws.Range("K132").Value = 11557.4121
$ws.Range("L132").Value = 7410.900000000001
$ws.Range("M132").Value = -9027.4121
$ws.Range("N132").Value = -12470.9
# row 136
$ws.Range("H136").Value = 1779.8684
$ws.Range("I136").Value = 791.7857
$ws.Range("K136").Value = 2375.3571
$ws.Range("M136").Value = 174.6428999999998
